$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28: fill in the new week's entry (2-3/6/2025)
$ws.Range("D28").Value = "2-3/6/2025"
$ws.Range("E28").Value = 127
$ws.Range("F28").Value = 234
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 1012
$ws.Range("J28").Value = "N/A"

# Row 29 & 30: set the date placeholders for upcoming weeks
$ws.Range("D29").Value = "4/6/2025"
$ws.Range("D30").Value = "5/6/2025"

# Update the view: scroll position and active cell selection
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("I29").Select()
